$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label for column B
$ws.Range("B1").Value = "Price per Sq. Ft. "

# Update existing price values
$ws.Range("B2").Value = 847
$ws.Range("B5").Value = 585

# Add two new rows of data
$ws.Range("A6").Value = "Philadelphia, PA"
$ws.Range("B6").Value = 127
$ws.Range("B6").NumberFormat = $ws.Range("B5").NumberFormat

$ws.Range("A7").Value = "Denver, CO"
$ws.Range("B7").Value = 289
$ws.Range("B7").NumberFormat = $ws.Range("B5").NumberFormat

# Move/center the active selection below the new data, matching the
# author's "Fix positioning to center" commit
[void]$ws.Range("B8").Select()
